$d = $word.ActiveDocument

$newBodyXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">Agenda </w:t></w:r><w:r><w:t>– Meeting minutes</w:t></w:r></w:p><w:p><w:r><w:t>1/12/23</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">Progress updates &lt; </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>5</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>mins</w:t></w:r></w:p><w:p><w:r><w:t>Nothing Really to show so far</w:t></w:r></w:p><w:p><w:r><w:t>Plan on working over the weekend</w:t></w:r></w:p><w:p/><w:p><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>What we learned from meeting with the client</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve"> (Monday follow up)</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve"> &lt; </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>30</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>mins</w:t></w:r></w:p><w:p><w:r><w:t>Phone app</w:t></w:r><w:r><w:t>:</w:t></w:r><w:r><w:t xml:space="preserve"> for tracking the order</w:t></w:r></w:p><w:p><w:r><w:t>Phone app</w:t></w:r><w:r><w:t>:</w:t></w:r><w:r><w:t xml:space="preserve"> will be able to accept orders</w:t></w:r><w:r><w:t>, show the driver the route (map interface) and show what package goes where</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Desktop app: </w:t></w:r><w:r><w:t>manages both vehicles and order, and gives a route to the orders</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r></w:p><w:p><w:r><w:t>Research</w:t></w:r><w:r><w:t>:</w:t></w:r><w:r><w:t xml:space="preserve"> DOT vehicles</w:t></w:r><w:r><w:t xml:space="preserve"> to figure out volume </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>the</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> can carry</w:t></w:r></w:p><w:p><w:r><w:t>Use Google maps to get route data</w:t></w:r></w:p><w:p><w:r><w:t>Security – google oAuth</w:t></w:r></w:p><w:p><w:r><w:t>Optimization: driver should be hauling both ways not driving with an empty truck if possible</w:t></w:r><w:r><w:t xml:space="preserve"> and batch orders</w:t></w:r></w:p><w:p><w:r><w:t>Language: React Native Xpo</w:t></w:r><w:r><w:t xml:space="preserve"> or </w:t></w:r><w:r><w:t>React Native Js</w:t></w:r><w:r><w:t xml:space="preserve"> or Node js</w:t></w:r><w:r><w:t xml:space="preserve"> (Because this is what the client wants)</w:t></w:r></w:p><w:p><w:r><w:t>APIs: Google Maps Api, security api</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r></w:p><w:p><w:r><w:t>Backend: mongo db or sql</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r></w:p><w:p><w:r><w:t>Quarkus or springboot?</w:t></w:r></w:p><w:p/><w:p><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Expected difficulties &lt; 5 mins</w:t></w:r></w:p><w:p><w:r><w:t>SRS might not be up the client’s specifications because we couldn’t hear him for half of the meeting</w:t></w:r></w:p><w:p/><w:p><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">Plan until next meeting &lt; </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>5</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve"> mins</w:t></w:r></w:p><w:p><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>-</w:t></w:r><w:r><w:t>Write SRS for next meeti</w:t></w:r><w:r><w:t>ng</w:t></w:r><w:r><w:t xml:space="preserve"> (atleast a rough draft)</w:t></w:r></w:p><w:p><w:r><w:t>Research what Quarkus and SpringBoot are and how they fit into the project</w:t></w:r></w:p><w:p><w:r><w:t>Research what backend we want to use and why</w:t></w:r></w:p><w:p><w:r><w:t>Research what language we want to use and why</w:t></w:r></w:p><w:p><w:r><w:lastRenderedPageBreak/><w:t>Write down any questions we come across when w</w:t></w:r><w:r><w:t>riting the SRS</w:t></w:r></w:p><w:p/><w:p/><w:p/><w:p/><w:p/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$d.Content.InsertXML($newBodyXml)

Write-Output ("Paragraph count: " + $d.Paragraphs.Count)
